$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value = [double]"1.145354048293079E-05"
$ws.Cells.Item(2, 3).Value = [double]"0.0002816624389738731"
$ws.Cells.Item(2, 4).Value = [double]"7.951149563867332E-06"
$ws.Cells.Item(2, 5).Value = [double]"3.67921254418313E-06"
$ws.Cells.Item(2, 6).Value = [double]"3.84664969921289E-06"
$ws.Cells.Item(2, 7).Value = [double]"6.995991586045136E-08"
$ws.Cells.Item(2, 8).Value = [double]"7.897454239236268E-06"
$ws.Cells.Item(2, 9).Value = [double]"3.083063053743729E-07"
$ws.Cells.Item(3, 2).Value = [double]"0.0006333609825530662"
$ws.Cells.Item(3, 3).Value = [double]"0.0001753018583174892"
$ws.Cells.Item(3, 4).Value = [double]"1.97813312365902E-05"
$ws.Cells.Item(3, 5).Value = [double]"2.69398390417841E-06"
$ws.Cells.Item(3, 6).Value = [double]"4.946358789670313E-06"
$ws.Cells.Item(3, 7).Value = [double]"5.328705965368385E-08"
$ws.Cells.Item(3, 8).Value = [double]"9.551046070699343E-06"
$ws.Cells.Item(3, 9).Value = [double]"2.766339517457084E-07"
$ws.Cells.Item(4, 2).Value = [double]"1.114891323972245E-05"
$ws.Cells.Item(4, 3).Value = [double]"0.0002905397748014238"
$ws.Cells.Item(4, 4).Value = [double]"7.46910221427477E-06"
$ws.Cells.Item(4, 5).Value = [double]"3.619328720061432E-06"
$ws.Cells.Item(4, 6).Value = [double]"3.738430381039493E-06"
$ws.Cells.Item(4, 7).Value = [double]"7.13178889277799E-08"
$ws.Cells.Item(4, 8).Value = [double]"7.797337738812971E-06"
$ws.Cells.Item(4, 9).Value = [double]"3.11740962821716E-07"
$ws.Cells.Item(5, 2).Value = [double]"0.02059645052341619"
$ws.Cells.Item(5, 3).Value = [double]"0.0001252677480902964"
$ws.Cells.Item(5, 4).Value = [double]"1.045074495409359E-05"
$ws.Cells.Item(5, 5).Value = [double]"9.38282756179909E-06"
$ws.Cells.Item(5, 6).Value = [double]"4.360638136802208E-06"
$ws.Cells.Item(5, 7).Value = [double]"1.090728460884094E-07"
$ws.Cells.Item(5, 8).Value = [double]"7.79388080918153E-06"
$ws.Cells.Item(5, 9).Value = [double]"3.215595709704644E-07"
$ws.Cells.Item(6, 2).Value = [double]"1.429661666385584E-06"
$ws.Cells.Item(6, 3).Value = [double]"0.004764677455612487"
$ws.Cells.Item(6, 4).Value = [double]"3.481112960557198E-06"
$ws.Cells.Item(6, 5).Value = [double]"6.323707134603396E-06"
$ws.Cells.Item(6, 6).Value = [double]"7.19375618696283E-06"
$ws.Cells.Item(6, 7).Value = [double]"4.532758457465036E-08"
$ws.Cells.Item(6, 8).Value = [double]"5.809569839931861E-05"
$ws.Cells.Item(6, 9).Value = [double]"2.772182927657241E-09"
$ws.Cells.Item(7, 2).Value = [double]"4.216664573064577E-09"
$ws.Cells.Item(7, 3).Value = [double]"0.0001947137330660098"
$ws.Cells.Item(7, 4).Value = [double]"8.846669779405403E-06"
$ws.Cells.Item(7, 5).Value = [double]"4.38156400450751E-06"
$ws.Cells.Item(7, 6).Value = [double]"6.792433822955246E-06"
$ws.Cells.Item(7, 7).Value = [double]"6.345664749218827E-08"
$ws.Cells.Item(7, 8).Value = [double]"9.703087466617077E-06"
$ws.Cells.Item(7, 9).Value = [double]"3.498223013487234E-07"
$ws.Cells.Item(8, 2).Value = [double]"2.197108210477694E-06"
$ws.Cells.Item(8, 3).Value = [double]"0.0001980130164593863"
$ws.Cells.Item(8, 4).Value = [double]"1.189320252614756E-05"
$ws.Cells.Item(8, 5).Value = [double]"4.590822228811575E-06"
$ws.Cells.Item(8, 6).Value = [double]"4.814066743019652E-06"
$ws.Cells.Item(8, 7).Value = [double]"5.756053234051706E-08"
$ws.Cells.Item(8, 8).Value = [double]"1.001396999433801E-05"
$ws.Cells.Item(8, 9).Value = [double]"2.57681175064306E-07"
$ws.Cells.Item(9, 2).Value = [double]"0.0003410168527335777"
$ws.Cells.Item(9, 3).Value = [double]"0.0002999616353533963"
$ws.Cells.Item(9, 4).Value = [double]"1.067804402676842E-05"
$ws.Cells.Item(9, 5).Value = [double]"4.368371011587158E-06"
$ws.Cells.Item(9, 6).Value = [double]"4.745522187876716E-06"
$ws.Cells.Item(9, 7).Value = [double]"5.806658842644197E-08"
$ws.Cells.Item(9, 8).Value = [double]"8.821391913160781E-06"
$ws.Cells.Item(9, 9).Value = [double]"2.905696052657447E-07"
$ws.Cells.Item(10, 2).Value = [double]"0.6499762046036415"
$ws.Cells.Item(10, 3).Value = [double]"0.4297917795193822"
$ws.Cells.Item(10, 4).Value = [double]"0.9856327348448009"
$ws.Cells.Item(10, 5).Value = [double]"0.9050345930354642"
$ws.Cells.Item(10, 6).Value = [double]"0.03362187756545043"
$ws.Cells.Item(10, 7).Value = [double]"0.0007363391154457943"
$ws.Cells.Item(10, 8).Value = [double]"0.0003857997176008876"
$ws.Cells.Item(10, 9).Value = [double]"2.855354205744121E-06"
